$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C column ("C02") values were stored three orders of magnitude too
# large (e.g. 9387.01... instead of 9.38701...) - rescale every data row
# (2-65) by dividing by 1000, keeping full floating point precision.
for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 / 1000
}

# Row 1 had an explicit custom height (17.25) left over from manual
# resizing; auto-fit it back so it reverts to the sheet's standard height
# with no explicit override.
$ws.Rows("1").AutoFit() | Out-Null

# Update the active selection left over from editing to reflect where the
# author last left the cursor.
$ws.Range("L24").Select() | Out-Null
